$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The video file used for the fiber-optic "sedan" listings was renamed from
# "fiber-optic.MOV" to "fiber-optic.sedan.mp4". Update both rows that
# reference it (row 4 = fiber-sedan, row 6 = fiber-suv).
$ws.Range("E4").Value = "fiber-optic.sedan.mp4"
$ws.Range("E6").Value = "fiber-optic.sedan.mp4"

# Move the selection to E6, matching where the edit was made, and scroll the
# sheet so row 5 is the first visible row (best-effort; some hosts only
# persist the active-cell selection, not the scroll/topLeftCell position).
$ws.Range("E6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
